$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row (rows 2-261).
# All of them change from 45182 (2023-09-13) to 45184 (2023-09-15).
$ws.Range("C2:C261").Value = 45184
